# Generate Report for Handback
# Updates the handoff/handback timestamps for the row corresponding to the
# "6cf6f32d-0ba7-4241-9914-7de703426469.md" file (row 3 on every sheet) to
# reflect the newly generated xliff report.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for row 3.
$overview.Range("G3").Value = "2016-10-24 09:26:28"

# zh-cn sheet: "Correspond Handoff Datetime" (H) / "Correspond Handback DateTime" (K) for row 3.
$zhcn.Range("H3").Value = "2016-10-24 09:26:17"
$zhcn.Range("K3").Value = "2016-10-24 09:27:08"

# de-de sheet: "Correspond Handoff Datetime" (H) / "Correspond Handback DateTime" (K) for row 3.
$dede.Range("H3").Value = "2016-10-24 09:26:28"
$dede.Range("K3").Value = "2016-10-24 09:27:25"
